# Ready to fit incidence. OK but not perfect.
# Update the incidence2018_plus sheet so that probabilities are computed by
# halving the prevalence2018 figures (instead of dividing by 100), for rows
# 27 through 122 of column C, then move the active selection/tab as the
# author left them.

$wb = $excel.ActiveWorkbook

$wsPrev = $wb.Worksheets.Item("prevalence2018")
$wsInc  = $wb.Worksheets.Item("incidence2018_plus")

# --- 1. Re-point the formulas in incidence2018_plus!C27:C122 -------------
# Was "=prevalence2018!C<row>/100", now "=prevalence2018!C<row>/2"
for ($r = 27; $r -le 122; $r++) {
    $wsInc.Range("C$r").Formula = "=prevalence2018!C$r/2"
}

# --- 2. Restore the selections left on each sheet -------------------------
$wsPrev.Activate()
$wsPrev.Range("C1:C1048576").Select()

$wsInc.Activate()
$wsInc.Range("F16").Select()

Write-Host "done"
